# "Restructure directories, add initializer for workbooks and worksheets"
#
# - "Main"   -> renamed to "Accounts"         (keep only the header row)
# - "Wealth" -> renamed to "Wealth Allocation" (keep only the header row)
# - a brand-new blank "Sheet1" is inserted in front of everything else and
#   becomes the initializer / active sheet for the workbook.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

function Clear-DataRowsKeepHeader($ws) {
    $lastRow = $ws.UsedRange.Rows.Count
    if ($lastRow -gt 1) {
        $ws.Range("A2:A" + $lastRow).EntireRow.Delete()
    }
}

# --- Accounts (was "Main"): Bank Name | Current Balance | Asset Type ---
$wsAccounts = $wb.Worksheets.Item("Main")
Clear-DataRowsKeepHeader $wsAccounts
$wsAccounts.Name = "Accounts"

# --- Wealth Allocation (was "Wealth"): Class | Balance ---
$wsWealth = $wb.Worksheets.Item("Wealth")
Clear-DataRowsKeepHeader $wsWealth
$wsWealth.Name = "Wealth Allocation"

# --- new blank initializer sheet, placed first and made active ---
$wsInit = $wb.Worksheets.Add($wsAccounts)
$wsInit.Name = "Sheet1"
$wsInit.Activate()
